$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) column cells whose new values look numeric,
# so Excel stores them as text (matching the source inlineStr cells) instead of
# auto-converting to numbers and losing formatting (e.g. trailing zeros).
$dCells = @("D2","D3","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D38","D39","D40","D41","D42","D44","D45","D46","D47","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = "64.517.23"
$ws.Range("E2").Value = "  -2.67%  "
$ws.Range("D3").Value = "3.362.53"
$ws.Range("E3").Value = "  -4.55%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "556.08"
$ws.Range("E5").Value = "  -4.82%  "
$ws.Range("D6").Value = "176.25"
$ws.Range("E6").Value = "  -1.61%  "
$ws.Range("D7").Value = "0.616"
$ws.Range("E7").Value = "  -2.37%  "
$ws.Range("D8").Value = "3.353.59"
$ws.Range("E8").Value = "  -4.62%  "
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.16%  "
$ws.Range("D10").Value = "0.627"
$ws.Range("E10").Value = "  -1.95%  "
$ws.Range("D11").Value = "0.162"
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("D12").Value = "54.20"
$ws.Range("E12").Value = "  -3.55%  "
$ws.Range("D13").Value = "0.0000273"
$ws.Range("E13").Value = "  -2.22%  "
$ws.Range("D14").Value = "9.07"
$ws.Range("E14").Value = "  -2.62%  "
$ws.Range("D15").Value = "3.888.63"
$ws.Range("E15").Value = "  -4.87%  "
$ws.Range("D16").Value = "18.40"
$ws.Range("E16").Value = "  -0.27%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "0.118"
$ws.Range("E17").Value = "  -2.89%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.355.27"
$ws.Range("E18").Value = "  -4.72%  "
$ws.Range("D19").Value = "11.85"
$ws.Range("E19").Value = "  -2.10%  "
$ws.Range("D20").Value = "64.403.05"
$ws.Range("E20").Value = "  -2.82%  "
$ws.Range("D21").Value = "0.982"
$ws.Range("E21").Value = "  -2.77%  "
$ws.Range("D22").Value = "428.11"
$ws.Range("E22").Value = "  +2.84%  "
$ws.Range("D23").Value = "4.89"
$ws.Range("E23").Value = "  +10.07%  "
$ws.Range("D24").Value = "4.11"
$ws.Range("E24").Value = "  -4.66%  "
$ws.Range("D25").Value = "84.18"
$ws.Range("E25").Value = "  -1.79%  "
$ws.Range("D26").Value = "13.33"
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("D27").Value = "10.78"
$ws.Range("E27").Value = "  -3.25%  "
$ws.Range("D28").Value = "2.84"
$ws.Range("E28").Value = "  -0.96%  "
$ws.Range("D29").Value = "8.74"
$ws.Range("E29").Value = "  -4.19%  "
$ws.Range("D30").Value = "29.80"
$ws.Range("E30").Value = "  -2.20%  "
$ws.Range("D31").Value = "6.68"
$ws.Range("E31").Value = "  +1.24%  "
$ws.Range("D32").Value = "11.47"
$ws.Range("E32").Value = "  -2.53%  "
$ws.Range("D33").Value = "573.41"
$ws.Range("E33").Value = "  -5.47%  "
$ws.Range("D34").Value = "0.108"
$ws.Range("E34").Value = "  -3.22%  "
$ws.Range("D35").Value = "58.29"
$ws.Range("E35").Value = "  -2.61%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("E37").Value = "  -7.59%  "
$ws.Range("D38").Value = "3.54"
$ws.Range("E38").Value = "  -3.81%  "
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "0.0₃0760"
$ws.Range("E39").Value = "  -5.70%  "
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").Value = "35.76"
$ws.Range("E40").Value = "  -4.22%  "
$ws.Range("D41").Value = "0.368"
$ws.Range("E41").Value = "  -4.44%  "
$ws.Range("D42").Value = "3.113.81"
$ws.Range("E42").Value = "  -4.78%  "
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").Value = "2.82"
$ws.Range("E44").Value = "  -5.69%  "
$ws.Range("D45").Value = "3.25"
$ws.Range("E45").Value = "  -4.11%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "0.0409"
$ws.Range("E46").Value = "  -3.02%  "
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").Value = "2.47"
$ws.Range("E47").Value = "  -3.19%  "
$ws.Range("E48").Value = "  -2.30%  "
$ws.Range("E49").Value = "  -4.74%  "
$ws.Range("D50").Value = "135.05"
$ws.Range("E50").Value = "  -3.45%  "
$ws.Range("D51").Value = "8.26"
$ws.Range("E51").Value = "  -5.08%  "
